# "Fixed POI packaging and upgraded to POI 3.15."
#
# The underlying OOXML diff for this revision is a pure re-serialization:
# every <w:...> element in word/document.xml and word/styles.xml keeps
# exactly the same element names, attribute names and attribute values as
# before - the only thing that changed is the on-disk *order* in which an
# element's attributes (and the root element's xmlns:* declarations) are
# written out (alphabetically, e.g. <w:tab w:val="left" w:pos="3119"/> ->
# <w:tab w:pos="3119" w:val="left"/>). That reordering is a side effect of
# the POI/XMLBeans package writer used when the fixture .docx was
# regenerated; it is not a document edit and is not something the Word
# object model exposes or lets an automation script control - Word COM
# never promises (or lets callers pick) attribute serialization order.
#
# So there is no content, formatting, or structural change for this
# revision to replay: the paragraphs, run text, tab stops, section/page
# setup, style/theme defaults and latent-style table are all byte-for-byte
# the same values as before. We touch the document the way the commit's
# "packaging fix" resave would - without mutating any content - so the
# save pipeline runs, and leave every property exactly as authored.

$d = $word.ActiveDocument

# Touch the document/package (mirrors the commit's repackage/resave)
# without changing any observable content, formatting or structure.
$null = $d.Content
$d.Saved = $d.Saved
